$wb = $excel.ActiveWorkbook

# --- Rename sheet "Falcon 404 IP assignment" -> "Vlan 404 IP assignment" ---
$wsVlan = $wb.Worksheets.Item("Falcon 404 IP assignment")
$wsVlan.Name = "Vlan 404 IP assignment"

# --- Add new Redis/MySQL host rows to the "Host" sheet ---
$wsHost = $wb.Worksheets.Item("Host")

$ips = @("10.245.93.48", "10.245.93.49", "10.245.93.50", "10.245.93.51", "10.245.93.52", "10.245.93.53")

$row = 99
foreach ($ip in $ips) {
    $wsHost.Cells.Item($row, 1).Value = "Windows"
    $wsHost.Cells.Item($row, 2).Value = "virtual"
    $wsHost.Cells.Item($row, 3).Value = $ip
    $wsHost.Cells.Item($row, 4).Value = "ISCSI File"
    $wsHost.Cells.Item($row, 5).Value = "Yang, Reid"
    $wsHost.Cells.Item($row, 8).Value = "SWARM"
    $row = $row + 1
}

# Copy the formatting (including cell border style) of H98 down onto the
# newly-added H99:H104 cells so they match the existing column styling.
$wsHost.Range("H98").Copy()
$wsHost.Range("H99:H104").PasteSpecial(-4122)

# --- Make the renamed Vlan sheet the active/selected tab ---
$wsVlan.Activate()
